$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D2:D51 and E2:E51 to be treated as text so that
# numeric-looking values (e.g. "2.60") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "61.519.32"
$ws.Range("E2").Value = "  -2.33%  "
$ws.Range("D3").Value = "2.948.70"
$ws.Range("E3").Value = "  -3.31%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "580.87"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "140.68"
$ws.Range("E6").Value = "  -7.21%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "0.519"
$ws.Range("E8").Value = "  -2.83%  "
$ws.Range("D9").Value = "2.946.43"
$ws.Range("E9").Value = "  -3.33%  "
$ws.Range("E10").Value = "  -6.04%  "
$ws.Range("E11").Value = "  -2.70%  "
$ws.Range("D12").Value = "0.456"
$ws.Range("E12").Value = "  +1.77%  "
$ws.Range("E13").Value = "  -4.21%  "
$ws.Range("D14").Value = "33.83"
$ws.Range("E14").Value = "  -6.42%  "
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("D16").Value = "3.437.14"
$ws.Range("E16").Value = "  -3.24%  "
$ws.Range("D17").Value = "6.95"
$ws.Range("E17").Value = "  -2.47%  "
$ws.Range("D18").Value = "61.543.94"
$ws.Range("E18").Value = "  -2.28%  "
$ws.Range("D19").Value = "2.947.65"
$ws.Range("E19").Value = "  -3.39%  "
$ws.Range("D20").Value = "448.07"
$ws.Range("E20").Value = "  -6.54%  "
$ws.Range("D21").Value = "13.77"
$ws.Range("E21").Value = "  -3.52%  "
$ws.Range("E22").Value = "  -4.03%  "
$ws.Range("E23").Value = "  -3.24%  "
$ws.Range("D24").Value = "80.92"
$ws.Range("E24").Value = "  -1.15%  "
$ws.Range("D25").Value = "12.04"
$ws.Range("E25").Value = "  -4.55%  "
$ws.Range("E26").Value = "  -10.85%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "9.47"
$ws.Range("E28").Value = "  -9.91%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "2.60"
$ws.Range("E30").Value = "  -2.27%  "
$ws.Range("D31").Value = "6.81"
$ws.Range("E31").Value = "  -7.81%  "
$ws.Range("E32").Value = "  -6.86%  "
$ws.Range("D33").Value = "27.03"
$ws.Range("E33").Value = "  -2.01%  "
$ws.Range("E34").Value = "  -4.14%  "
$ws.Range("E35").Value = "  -5.36%  "
$ws.Range("D36").Value = "0.0₃0770"
$ws.Range("E36").Value = "  -5.57%  "
$ws.Range("E37").Value = "  -4.56%  "
$ws.Range("D38").Value = "2.06"
$ws.Range("E38").Value = "  -6.53%  "
$ws.Range("D39").Value = "49.96"
$ws.Range("E39").Value = "  -0.76%  "
$ws.Range("D40").Value = "9.06"
$ws.Range("E40").Value = "  -1.75%  "
$ws.Range("D41").Value = "0.118"
$ws.Range("E41").Value = "  +2.90%  "
$ws.Range("E42").Value = "  -14.09%  "
$ws.Range("D43").Value = "386.23"
$ws.Range("E43").Value = "  -10.18%  "
$ws.Range("D44").Value = "0.0350"
$ws.Range("E44").Value = "  -3.22%  "
$ws.Range("D45").Value = "2.704.96"
$ws.Range("E45").Value = "  -4.56%  "
$ws.Range("E46").Value = "  -8.99%  "
$ws.Range("D47").Value = "36.60"
$ws.Range("E47").Value = "  -4.10%  "
$ws.Range("D48").Value = "129.66"
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("E51").Value = "  -2.38%  "

# Remove the temporary text number-format so the cells return to the
# workbook default style (no explicit style index), matching the source.
$ws.Range("D2:D51").ClearFormats()
$ws.Range("E2:E51").ClearFormats()
